$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (B1:E1)
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 data values (B2:E2)
$ws.Range("B2").Value = 20.273753141035748
$ws.Range("C2").Value = -3.408208246140525
$ws.Range("D2").Value = -1.7819954847074513
$ws.Range("E2").Value = 1.9413154833607607

# Row 3 data values (B3:E3)
$ws.Range("B3").Value = 18.467710164024432
$ws.Range("C3").Value = 2.838265024556307
$ws.Range("D3").Value = -12.441338495981078
$ws.Range("E3").Value = 11.554444518804551

# Update selection to match the new active range
$ws.Range("B1:E3").Select()
